$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values (swap with row 3's original values)
$ws.Range("D2").Value = 44273
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 14000
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 233

# Row 3 values (swap with row 2's original values)
$ws.Range("D3").Value = 44291
$ws.Range("J3").Value = 25
$ws.Range("K3").Value = 11000
$ws.Range("L3").Value = 11000
$ws.Range("M3").Value = 11000
$ws.Range("O3").Value = "Limache"
$ws.Range("P3").Value = 183
